$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (e.g. "25.876.61" / "210.10"); force Text
# format before writing so Excel does not auto-convert numeric-looking
# strings ("210.10", "0.481", ...) into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.876.61'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.585.87'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '210.10'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.481'
$ws.Range('E7').Value = '  -3.40%  '
$ws.Range('E8').Value = '  -0.85%  '
$ws.Range('D9').Value = '0.0617'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = '18.08'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').Value = '0.0791'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.805.67'
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('D13').Value = '1.587.87'
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('E15').Value = '  -2.81%  '
$ws.Range('D16').Value = '25.868.73'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '0.0₃0724'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').Value = '59.92'
$ws.Range('E18').Value = '  -2.80%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '191.50'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').Value = '5.93'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('D24').Value = '0.133'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').Value = '142.03'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = '15.11'
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('D29').Value = '6.45'
$ws.Range('E29').Value = '  -2.89%  '
$ws.Range('E30').Value = '  -5.33%  '
$ws.Range('D31').Value = '0.0471'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').Value = '3.03'
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').Value = '1.103.05'
$ws.Range('E36').Value = '  -1.91%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '2.34'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '0.505'
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.829'
$ws.Range('E40').Value = '  +10.62%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0150'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('D42').Value = '0.780'
$ws.Range('E42').Value = '  -7.54%  '
$ws.Range('D43').Value = '5.19'
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').Value = '93.86'
$ws.Range('E44').Value = '  -4.01%  '
$ws.Range('D45').Value = '1.719.20'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('D46').Value = '0.0₆0106'
$ws.Range('E46').Value = '  -6.33%  '
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').Value = '53.31'
$ws.Range('E48').Value = '  -1.58%  '
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').Value = '0.406'
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('E51').Value = '  -0.14%  '
